$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 100 - this shifts the former rows 100-102 down to 101-103
$ws.Rows(100).Insert()

# Populate the newly inserted row 100 with the new weekly record
$ws.Cells.Item(100, 1).Value = 3
$ws.Cells.Item(100, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(100, 3).Value = "Coquimbo"
$ws.Cells.Item(100, 4).Value = 45041
$ws.Cells.Item(100, 5).Value = 5
$ws.Cells.Item(100, 6).Value = 100112035
$ws.Cells.Item(100, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(100, 8).Value = "Sin especificar"
$ws.Cells.Item(100, 9).Value = "Primera"
$ws.Cells.Item(100, 10).Value = 65
$ws.Cells.Item(100, 11).Value = 12000
$ws.Cells.Item(100, 12).Value = 12000
$ws.Cells.Item(100, 13).Value = 12000
$ws.Cells.Item(100, 14).Value = "`$/malla 10 kilos"
$ws.Cells.Item(100, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(100, 16).Value = 1200
$ws.Cells.Item(100, 17).Value = 10
$ws.Cells.Item(100, 18).Value = "Hortaliza"
